$wb = $excel.ActiveWorkbook

# Both "展览" (Exhibition) and "全部类型" (All Types) sheets carry the same
# duplicated data; the diff bumps the "想去人数" (want-to-go count) values
# in F3 (96 -> 97) and F9 (368 -> 369) on each of them.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 97
    $ws.Range("F9").Value = 369
}
